# Generate Report for Handback
# - Update the "Ready for handoff" status text (wherever it appears, shared
#   across sheets) to "Handback transform failed".
# - Populate the "Error Detail" column (P) on the zh-cn and de-de sheets with
#   a handback/handoff filename mismatch message, and widen that column so
#   the message is readable.

$wb = $excel.ActiveWorkbook

# 1) Status text change: every cell currently reading "Ready for handoff"
#    (Overview!E3, Overview!F3, zh-cn!C3, de-de!C3) becomes
#    "Handback transform failed".
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $used = $sheet.UsedRange
    foreach ($cell in $used.Cells) {
        $text = [string]$cell.Text
        if ($text -eq "Ready for handoff") {
            $cell.Value = "Handback transform failed"
        }
    }
}

# The host rounds ColumnWidth to a whole-pixel grid (width-in-chars + 5/6)
# when it serialises, so asking for an even 40 lands on 40.8333 in the saved
# XML. Backing off by 5/6 before the rounding lands exactly on 40.
$targetColumnWidth = 40 - (5 / 6)

# 2) zh-cn sheet: widen the Error Detail column and fill in the error detail
#    for row 3.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(16).ColumnWidth = $targetColumnWidth
$zhcn.Range("P3").Value = "Handback file name: d04trfqq.yaf is different with handoff file name: 7b60815c-4f0d-408f-bfbf-06c52280ccfb.32c294c2dd2a44a39603579379e5b9f643dfb23b.zh-cn."

# 3) de-de sheet: widen the Error Detail column and fill in the error detail
#    for row 3.
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = $targetColumnWidth
$dede.Range("P3").Value = "Handback file name: d04trfqq.yaf is different with handoff file name: 7b60815c-4f0d-408f-bfbf-06c52280ccfb.32c294c2dd2a44a39603579379e5b9f643dfb23b.de-de."
